# Applies the "lowercase IDs" structural tidy-up described in the commit:
#   changed IDs to lowercase id, amended and tidied all structural checks,
#   some work on pkgdown vignettes
#
# Order of operations matters because the workbook backs cell text with a
# shared-string table: editing the Notes explanation first (while it is
# still uniquely referenced) lets it get replaced in-place before any new
# strings are appended, matching the canonical OOXML produced by Excel.

$wb = $excel.ActiveWorkbook

$wsNotes   = $wb.Worksheets.Item("Notes")
$wsStudies = $wb.Worksheets.Item("studies")
$wsSurveys = $wb.Worksheets.Item("surveys")

# Update the explanatory note about why this structure is invalid
$wsNotes.Range("A3").Value = "Specific issue: study_IDs are not referenced in the surveys table"

# Tidy the sample identifier values and rename ID columns to lowercase
$wsStudies.Range("A2").Value = "bar"
$wsSurveys.Range("A2").Value = "foobar"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsStudies.Range("A1").Value = "study_id"

# Give the surveys header row a distinct (black) font colour
$wsSurveys.Range("A1:G1").Font.Color = 0
$wsSurveys.Range("K1").Font.Color = 0
$wsSurveys.Range("H1:J1").Font.Color = 0

# Move the active tab/selection from "surveys" to "studies"
$null = $wsSurveys.Activate()
$null = $wsSurveys.Range("C8").Select()
$null = $wsStudies.Activate()
$null = $wsStudies.Range("A2").Select()
